# Registers.xlsx edit: add DAC80508_CONFIG sheet, fix DAC80508 typos/values
$wb = $excel.ActiveWorkbook
$ws6 = $wb.Worksheets.Item("DAC80508")

# --- Fix up the existing DAC80508 sheet ---
# "CONFI" -> "CONFIG" (row 5, column A)
$ws6.Range("A5").Value = "CONFIG"

# widen column A slightly to fit the longer labels
$ws6.Columns.Item(1).ColumnWidth = 13.333333333333332

# --- Add the new DAC80508_CONFIG worksheet after DAC80508 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws7.Name = "DAC80508_CONFIG"

$ws7.Columns.Item(1).ColumnWidth = 13.0
$ws7.Columns.Item(2).ColumnWidth = 11.833333333333332
$ws7.Columns.Item(3).ColumnWidth = 13.833333333333332

$ws7.Range("A1").Value = "Name"
$ws7.Range("B1").Value = "Hex Address"
$ws7.Range("C1").Value = "Default Value"
$ws7.Range("D1").Value = "Bit Width"
$ws7.Range("E1").Value = "Bit Index"

$names = @("ALM-SEL","ALM-EN","CRC-EN","FSDO","DSDO","REF-PWDWN","DAC7-PWDWN","DAC6-PWDWN","DAC5-PWDWN","DAC4-PWDWN","DAC3-PWDWN","DAC2-PWDWN","DAC1-PWDWN","DAC0-PWDWN")

$row = 2
$bitIndex = 13
foreach ($n in $names) {
    $ws7.Cells.Item($row, 1).Value = $n
    $ws7.Cells.Item($row, 2).Value = "0x3"
    $ws7.Cells.Item($row, 3).Value = "0x0"
    $ws7.Cells.Item($row, 4).Value = 1
    $ws7.Cells.Item($row, 5).Value = $bitIndex
    $row++
    $bitIndex--
}

$ws7.Range("A1:E1").Select()

# --- Finish fixing up DAC80508 sheet: DEVICE ID -> ID, and its default value ---
$ws6.Range("A3").Value = "ID"
$ws6.Range("C3").Value = "0x0896"

# Restore DAC80508 as the active sheet/selection
$ws6.Activate()
$ws6.Range("C4").Select()
